# Updated cryptos list on Fri Mar 15 15:28:37 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ D='68.365.29'; E='  -4.26%  ' }
    3 = @{ D='3.700.04'; E='  -4.39%  ' }
    4 = @{ E='  -0.16%  ' }
    5 = @{ D='585.00'; E='  -2.17%  ' }
    6 = @{ D='183.04'; E='  +9.16%  ' }
    7 = @{ D='3.690.33'; E='  -4.58%  ' }
    8 = @{ D='0.628'; E='  -6.97%  ' }
    9 = @{ D='1.00'; E='  +0.13%  ' }
    10 = @{ D='0.718'; E='  -4.85%  ' }
    11 = @{ D='0.164'; E='  -7.67%  ' }
    12 = @{ D='55.66'; E='  +4.27%  ' }
    13 = @{ E='  -9.40%  ' }
    14 = @{ D='10.40'; E='  -7.06%  ' }
    15 = @{ D='4.200.79'; E='  -6.90%  ' }
    16 = @{ D='3.709.97'; E='  -4.92%  ' }
    17 = @{ D='19.40'; E='  -6.71%  ' }
    18 = @{ E='  -2.56%  ' }
    19 = @{ D='1.13'; E='  -7.01%  ' }
    20 = @{ D='12.76'; E='  -8.03%  ' }
    21 = @{ D='68.115.96'; E='  -4.67%  ' }
    22 = @{ D='409.18'; E='  -5.82%  ' }
    23 = @{ D='4.48'; E='  -5.66%  ' }
    24 = @{ D='88.56'; E='  -6.25%  ' }
    25 = @{ E='  -7.89%  ' }
    26 = @{ D='11.08'; E='  +1.62%  ' }
    27 = @{ D='12.80'; E='  -7.30%  ' }
    28 = @{ E='  -3.65%  ' }
    29 = @{ D='6.05'; E='  +1.86%  ' }
    30 = @{ D='9.49'; E='  -6.76%  ' }
    31 = @{ D='32.69'; E='  -6.64%  ' }
    32 = @{ E='  -5.10%  ' }
    33 = @{ D='12.53'; E='  -7.54%  ' }
    34 = @{ B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.117'; E='  -6.52%  ' }
    35 = @{ B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='65.45'; E='  -5.45%  ' }
    36 = @{ D='43.62'; E='  -13.66%  ' }
    37 = @{ D='596.55'; E='  -3.35%  ' }
    38 = @{ D='0.0₃0891'; E='  -9.75%  ' }
    39 = @{ B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.402'; E='  -4.00%  ' }
    40 = @{ B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.00'; E='  +0.07%  ' }
    41 = @{ E='  -0.13%  ' }
    42 = @{ E='  -3.93%  ' }
    43 = @{ E='  +5.12%  ' }
    44 = @{ E='  -8.43%  ' }
    45 = @{ E='  -7.42%  ' }
    46 = @{ E='  -7.02%  ' }
    47 = @{ D='9.30'; E='  -8.90%  ' }
    48 = @{ D='2.773.75'; E='  -2.36%  ' }
    49 = @{ B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.134'; E='  -6.93%  ' }
    50 = @{ B='WEMIXToken'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='2.69'; E='  -2.31%  ' }
    51 = @{ B='ApeXProtocol'; C='https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'; D='3.15'; E='  -6.75%  ' }
}

foreach ($r in $changes.Keys) {
    $row = $changes[$r]
    foreach ($col in @('B', 'C', 'D', 'E')) {
        if ($row.ContainsKey($col)) {
            $value = $row[$col]
            $cellRef = "$col$r"
            if ($value -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
                # Value would otherwise be auto-coerced to a number by Excel;
                # force text storage (matches source workbook's inline-string cells)
                # while keeping the General number format on the cell.
                $ws.Range($cellRef).NumberFormat = "@"
                $ws.Range($cellRef).Value = $value
                $ws.Range($cellRef).Style = "Normal"
            } else {
                $ws.Range($cellRef).Value = $value
            }
        }
    }
}
